$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.333.57"
$ws.Range("E2").Value = "  -3.04%  "

# Row 3
$ws.Range("D3").Value = "1.854.85"
$ws.Range("E3").Value = "  -3.10%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.34%  "

# Row 6
$ws.Range("E6").Value = "  +0.04%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4605"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.62%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3935"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.12%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.98"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -11.52%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07942"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.62%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.012"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.09%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.46"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.82%  "

# Row 13
$ws.Range("D13").Value = "1.868.33"
$ws.Range("E13").Value = "  -3.05%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.920"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.36%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.144"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.85%  "

# Row 16
$ws.Range("E16").Value = "  -0.05%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.91%  "

# Row 18
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001031"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.31%  "

# Row 19
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06576"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.50%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.54%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.21%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.463"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.63%  "

# Row 23
$ws.Range("D23").Value = "27.345.17"
$ws.Range("E23").Value = "  -3.07%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.07%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.295"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.28%  "

# Row 26
$ws.Range("D26").Value = "2.083.78"
$ws.Range("E26").Value = "  -3.52%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.37%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.37%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.066"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.28%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.451"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.11%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.36%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09415"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.21%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9492"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.33%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.442"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.10%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.585"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.59%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.265"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.79%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06037"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.96%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02227"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.93%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.206"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.67%  "

# Row 40
$ws.Range("E40").Value = "  +0.13%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.024"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.88%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5934"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.35%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1889"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.65%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.33%  "

# Row 45
$ws.Range("E45").Value = "  -1.56%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5628"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.65%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.05%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.397"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.86%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.916"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.14%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06755"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.99%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "109.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.67%  "
